# "Works by Year filter" -- adds a new "Migraine" worksheet (implied disease
# burden by race/ethnicity for severe headache/migraine), tweaks the
# selection on Expected_demographics, and makes Migraine the active sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new "Migraine" sheet after the last existing sheet (HCV).
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "Migraine"

# ---------------------------------------------------------------------
# 2. Populate the category labels first (column C, rows 3-9) so the
#    shared-string table fills in the same order as the source sheet.
# ---------------------------------------------------------------------
$ws3.Range("C3").Value = "White only"
$ws3.Range("C4").Value = "Black only"
$ws3.Range("C5").Value = "American Indian or Alaska Native only"
$ws3.Range("C6").Value = "Asian only"
$ws3.Range("C7").Value = "Native Hawaiian only"
$ws3.Range("C8").Value = "2 or more races"
$ws3.Range("C9").Value = "Hispanic or Latino"

# Headers
$ws3.Range("D2").Value = "% of adults with severe headache or migraine"
$ws3.Range("C1").Value = "https://www.cdc.gov/nchs/data/hus/2017/041.pdf"
$ws3.Range("E2").Value = "% of population"
$ws3.Range("F2").Value = "implied burden"
$ws3.Range("E11").Value = "total"

# ---------------------------------------------------------------------
# 3. Numeric inputs.
# ---------------------------------------------------------------------
$ws3.Range("D3").Value = 16
$ws3.Range("E3").Value = 76.5

$ws3.Range("D4").Value = 14.6
$ws3.Range("E4").Value = 14.3

$ws3.Range("D5").Value = 16.4

$ws3.Range("D6").Value = 9.9
$ws3.Range("E6").Value = 6.8

$ws3.Range("D8").Value = 22.5

$ws3.Range("D9").Value = 15.3
$ws3.Range("E9").Value = 18.5

# ---------------------------------------------------------------------
# 4. Formulas - implied burden (col F) and share of total (col G).
#    F4:F9 is filled as one range-formula (Excel encodes this as a
#    shared formula, si="0", anchored at F4) and the three rows with no
#    source data (5, 7, 8) are cleared back out again, leaving just an
#    empty, formatted cell - matching the source table's gaps.
# ---------------------------------------------------------------------
$ws3.Range("F3").Formula = "=(D3/100)*E3"
$ws3.Range("F4:F9").Formula = "=(D4/100)*E4"
$ws3.Range("F5").ClearContents()
$ws3.Range("F7").ClearContents()
$ws3.Range("F8").ClearContents()
$ws3.Range("F11").Formula = "=SUM(F3:F9)"

$ws3.Range("G3").Formula = "=F3/`$F`$11"
$ws3.Range("G4").Formula = "=F4/`$F`$11"
$ws3.Range("G6").Formula = "=F6/`$F`$11"
$ws3.Range("G9").Formula = "=F9/`$F`$11"

# ---------------------------------------------------------------------
# 5. Number formats. Column F = "0.0", column G = "0.0%" (the latter
#    reuses the workbook's existing Percent style). Also apply the "0.0"
#    format to the still-empty F/G cells on rows that have no data so
#    the styling matches the fully built-out table.
# ---------------------------------------------------------------------
$ws3.Range("F3:F9").NumberFormat = "0.0"
$ws3.Range("F11").NumberFormat = "0.0"

$ws3.Range("G3").NumberFormat = "0.0%"
$ws3.Range("G4").NumberFormat = "0.0%"
$ws3.Range("G5").NumberFormat = "0.0%"
$ws3.Range("G6").NumberFormat = "0.0%"
$ws3.Range("G7").NumberFormat = "0.0%"
$ws3.Range("G8").NumberFormat = "0.0%"
$ws3.Range("G9").NumberFormat = "0.0%"

# ---------------------------------------------------------------------
# 6. Column widths (best-fit like the source workbook).
# ---------------------------------------------------------------------
$ws3.Columns.Item(3).ColumnWidth = 31.998697916666668
$ws3.Columns.Item(4).ColumnWidth = 38.666666666666664
$ws3.Columns.Item(5).ColumnWidth = 13.166666666666666

# ---------------------------------------------------------------------
# 7. Update the selection on Expected_demographics (A5 -> B5:F5).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Expected_demographics")
$ws1.Select() | Out-Null
$ws1.Range("B5:F5").Select() | Out-Null

# ---------------------------------------------------------------------
# 8. Finish with Migraine as the active sheet / selection (matches the
#    "Works by Year filter" tab state in the target workbook).
# ---------------------------------------------------------------------
$ws3.Select() | Out-Null
$ws3.Range("G3").Select() | Out-Null
